$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 912.6
$ws.Range("J127").Value = 2216
$ws.Range("L127").Value = 6648
$ws.Range("N127").Value = -16568

$ws.Range("H132").Value = 11966692
$ws.Range("I132").Value = 14359640
$ws.Range("J132").Value = 1950
$ws.Range("K132").Value = 43078920
$ws.Range("L132").Value = 5850
$ws.Range("M132").Value = -43076390
$ws.Range("N132").Value = -10910

$ws.Range("H137").Value = 12565.634
$ws.Range("I137").Value = 23801.857
$ws.Range("J137").Value = 2733.9375
$ws.Range("K137").Value = 71405.571
$ws.Range("L137").Value = 8201.8125
$ws.Range("M137").Value = -68855.571
$ws.Range("N137").Value = -13301.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17684.861
$ws.Range("I32").Value = 18103.428
$ws.Range("K32").Value = 18103.428
$ws.Range("M32").Value = -17816.428

$ws.Range("H61").Value = 11673.363
$ws.Range("I61").Value = 1801.6666
$ws.Range("J61").Value = 23519.4
$ws.Range("K61").Value = 1801.6666
$ws.Range("L61").Value = 23519.4
$ws.Range("M61").Value = -1589.6666
$ws.Range("N61").Value = -23943.4

$ws.Range("H63").Value = 2259.3845
$ws.Range("I63").Value = 1708.4445
$ws.Range("K63").Value = 1708.4445
$ws.Range("M63").Value = -1022.4445

$ws.Range("H66").Value = 2259.3845
$ws.Range("I66").Value = 1708.4445
$ws.Range("K66").Value = 8542.2225
$ws.Range("M66").Value = -5110.2225

$ws.Range("H74").Value = 255141.2
$ws.Range("I74").Value = 300944.8
$ws.Range("K74").Value = 300944.8
$ws.Range("M74").Value = -300070.8

$ws.Range("H77").Value = 255141.2
$ws.Range("I77").Value = 300944.8
$ws.Range("K77").Value = 1504724
$ws.Range("M77").Value = -1500356

$ws.Range("H136").Value = 11673.363
$ws.Range("I136").Value = 1801.6666
$ws.Range("J136").Value = 23519.4
$ws.Range("K136").Value = 5404.9998
$ws.Range("L136").Value = 70558.20000000001
$ws.Range("M136").Value = -2854.9998
$ws.Range("N136").Value = -75658.20000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 42407.8
$ws.Range("J20").Value = 1508.5714
$ws.Range("L20").Value = 1508.5714
$ws.Range("N20").Value = -2002.5714

$ws.Range("H86").Value = 1726.1
$ws.Range("I86").Value = 1741.5333
$ws.Range("K86").Value = 1741.5333
$ws.Range("M86").Value = -618.5333000000001

$ws.Range("H89").Value = 1726.1
$ws.Range("I89").Value = 1741.5333
$ws.Range("K89").Value = 8707.666499999999
$ws.Range("M89").Value = -3091.666499999999

$ws.Range("H134").Value = 2366.795
$ws.Range("I134").Value = 1864.3667
$ws.Range("K134").Value = 5593.1001
$ws.Range("M134").Value = -3058.1001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 770.1818
$ws.Range("I22").Value = 295.85715
$ws.Range("J22").Value = 1600.25
$ws.Range("K22").Value = 295.85715
$ws.Range("L22").Value = 1600.25
$ws.Range("M22").Value = 54.14285000000001
$ws.Range("N22").Value = -2300.25

$ws.Range("H86").Value = 51938.223
$ws.Range("J86").Value = 21221.777
$ws.Range("L86").Value = 21221.777
$ws.Range("N86").Value = -23467.777

$ws.Range("H89").Value = 51938.223
$ws.Range("J89").Value = 21221.777
$ws.Range("L89").Value = 106108.885
$ws.Range("N89").Value = -117340.885

$ws.Range("H99").Value = 6802.6665
$ws.Range("I99").Value = 4118.3076
$ws.Range("K99").Value = 4118.3076
$ws.Range("M99").Value = -2620.3076

$ws.Range("H107").Value = 834.9474
$ws.Range("I107").Value = 644.9167
$ws.Range("J107").Value = 1160.7142
$ws.Range("K107").Value = 644.9167
$ws.Range("L107").Value = 1160.7142
$ws.Range("M107").Value = 1275.0833
$ws.Range("N107").Value = -5000.7142

$ws.Range("H126").Value = 6802.6665
$ws.Range("I126").Value = 4118.3076
$ws.Range("K126").Value = 12354.9228
$ws.Range("M126").Value = -9884.9228

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 398
$ws.Range("J7").Value = 430
$ws.Range("L7").Value = 1290
$ws.Range("N7").Value = -1514

$ws.Range("H131").Value = 165156.84
$ws.Range("J131").Value = 1981
$ws.Range("L131").Value = 5943
$ws.Range("N131").Value = -16023

$ws.Range("H132").Value = 1880
$ws.Range("I132").Value = 2045.6666
$ws.Range("J132").Value = 1134.5
$ws.Range("K132").Value = 18410.9994
$ws.Range("L132").Value = 10210.5
$ws.Range("M132").Value = -15880.9994
$ws.Range("N132").Value = -15270.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6958.8335
$ws.Range("I70").Value = 6699.12
$ws.Range("K70").Value = 6699.12
$ws.Range("M70").Value = -6429.12

$ws.Range("H73").Value = 6958.8335
$ws.Range("I73").Value = 6699.12
$ws.Range("K73").Value = 6699.12
$ws.Range("M73").Value = -5763.12

$ws.Range("H80").Value = 14203.8
$ws.Range("I80").Value = 3200
$ws.Range("J80").Value = 16954.75
$ws.Range("K80").Value = 3200
$ws.Range("L80").Value = 16954.75
$ws.Range("M80").Value = -2202
$ws.Range("N80").Value = -18950.75

$ws.Range("H83").Value = 14203.8
$ws.Range("I83").Value = 3200
$ws.Range("J83").Value = 16954.75
$ws.Range("K83").Value = 16000
$ws.Range("L83").Value = 84773.75
$ws.Range("M83").Value = -11008
$ws.Range("N83").Value = -94757.75

$ws.Range("H102").Value = 16557.537
$ws.Range("I102").Value = 23319.678
$ws.Range("K102").Value = 23319.678
$ws.Range("M102").Value = -21697.678

$ws.Range("H132").Value = 4307.7334
$ws.Range("I132").Value = 4258.2856
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 12774.8568
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -10244.8568
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 14000
$ws.Range("J44").Value = 14000
$ws.Range("L44").Value = 14000
$ws.Range("N44").Value = -14912

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = ""
$ws.Range("N100").Value = ""

$ws.Range("H122").Value = 5759.143
$ws.Range("I122").Value = 3470.5715
$ws.Range("K122").Value = 10411.7145
$ws.Range("M122").Value = -7961.7145

$ws.Range("H136").Value = 4170
$ws.Range("I136").Value = 4170
$ws.Range("K136").Value = 12510
$ws.Range("M136").Value = -9960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5065.6665
$ws.Range("J62").Value = 5498.5
$ws.Range("L62").Value = 5498.5
$ws.Range("N62").Value = -6746.5

$ws.Range("H65").Value = 5065.6665
$ws.Range("J65").Value = 5498.5
$ws.Range("L65").Value = 27492.5
$ws.Range("N65").Value = -33732.5

$ws.Range("H136").Value = 11312
$ws.Range("I136").Value = 12043.56
$ws.Range("J136").Value = 5215.6665
$ws.Range("K136").Value = 36130.68
$ws.Range("L136").Value = 15646.9995
$ws.Range("M136").Value = -33580.68
$ws.Range("N136").Value = -20746.9995
